# Update "想去人数" (interested-count) values in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets, matching the
# data refresh captured in the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row => new value for worksheet "展览"
$exhibitionUpdates = @{
    3  = 283
    4  = 1813
    10 = 3577
    13 = 74
    14 = 53
    17 = 115
    18 = 785
    19 = 8
    20 = 215
    21 = 134
    22 = 60
    23 = 68
    24 = 72
    25 = 2812
    26 = 5279
    29 = 481
    30 = 3102
    31 = 299
    32 = 2288
    34 = 494
    35 = 90
    36 = 142
    37 = 194
    39 = 45
    40 = 472
    41 = 817
    43 = 19
    46 = 501
}

# Row => new value for worksheet "全部类型"
$allTypesUpdates = @{
    3  = 283
    4  = 1813
    10 = 3577
    13 = 74
    15 = 53
    18 = 115
    19 = 785
    20 = 8
    21 = 215
    22 = 134
    23 = 60
    24 = 68
    25 = 72
    26 = 2812
    27 = 5279
    30 = 481
    31 = 3102
    32 = 299
    33 = 2288
    35 = 494
    36 = 90
    37 = 142
    38 = 194
    40 = 45
    41 = 472
    42 = 817
    44 = 19
    47 = 501
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
